$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 43, shifting rows 43:79 down to 44:80
$ws.Rows.Item(43).Insert()

# Populate the newly inserted row 43 with the new record
$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(43, 3).Value = "Metropolitana"
$ws.Cells.Item(43, 4).Value = 45062
$ws.Cells.Item(43, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(43, 5).Value = 13
$ws.Cells.Item(43, 6).Value = 100112035
$ws.Cells.Item(43, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(43, 8).Value = "Sin especificar"
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 52
$ws.Cells.Item(43, 11).Value = 22000
$ws.Cells.Item(43, 12).Value = 25000
$ws.Cells.Item(43, 13).Value = 23500
$ws.Cells.Item(43, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(43, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(43, 16).Value = 1567
$ws.Cells.Item(43, 17).Value = 15
$ws.Cells.Item(43, 18).Value = "Hortaliza"
